$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-44: update transaction date, product, quantity ---
$newDate = 46053.4375
$ws.Range("A2:A44").Value = $newDate

$rows = @(
    @{Row=2; Product='Excel Salmon Kuning'; Qty=1},
    @{Row=3; Product='excel tuna hijau'; Qty=5},
    @{Row=4; Product='Maxi 1kg'; Qty=3},
    @{Row=5; Product='Cat Choize Pink Kitten Tuna - 1 kg'; Qty=8},
    @{Row=6; Product='Life cat kaleng Tuna Kitten Pink - 400gr'; Qty=5},
    @{Row=7; Product='Pasir Taro 10 liter'; Qty=1},
    @{Row=8; Product='lezato salmon adult'; Qty=2},
    @{Row=9; Product='Life cat Kaleng tuna adult hijau 400 gr'; Qty=1},
    @{Row=10; Product='excel tuna ungu'; Qty=3},
    @{Row=11; Product='Bolt Mother Kitten tuna- 500g'; Qty=4},
    @{Row=12; Product='Furlove kitten salmon - 1kg'; Qty=2},
    @{Row=13; Product='Beauty premium Chicken and Salmon - 1kg'; Qty=1},
    @{Row=14; Product='Disposable syringe 3 cc/mL'; Qty=1},
    @{Row=15; Product='susu kitten growsy'; Qty=1},
    @{Row=16; Product='Excel mother kitten'; Qty=5},
    @{Row=17; Product='Cat Choize Kuning Kitten Salmon - 1 kg'; Qty=7},
    @{Row=18; Product='Markotop pouch Tuna Salmon'; Qty=2},
    @{Row=19; Product='Bolt Mother & Kitten Salmon- 500g'; Qty=2},
    @{Row=20; Product='Whiskas Pouch Junior Tuna Flavour'; Qty=2},
    @{Row=21; Product='Chester'; Qty=1},
    @{Row=22; Product='Pasir Repack 1 kg'; Qty=2},
    @{Row=23; Product='Bolt Donat tuna kuning 800 gr'; Qty=2},
    @{Row=24; Product='Bolt Tuna - 800g'; Qty=10},
    @{Row=25; Product='Cat Choize Hijau Tuna Adult - 800g'; Qty=4},
    @{Row=26; Product='Cat Choize Oren Salmon Adult - 800 gr'; Qty=2},
    @{Row=27; Product='Nice Tuna Adult 800 gr'; Qty=2},
    @{Row=28; Product='crystal kitty adult'; Qty=1},
    @{Row=29; Product='Snack Curah kucing stik all varian rasa'; Qty=5},
    @{Row=30; Product='Pasir Junior 5 L'; Qty=2},
    @{Row=31; Product='Susu Top Growth'; Qty=1},
    @{Row=32; Product='crystal kitty tuna chicken - mother kitten'; Qty=2},
    @{Row=33; Product='mister puss salmon - 500g'; Qty=1},
    @{Row=34; Product='Cat Choize Mother Kitten'; Qty=1},
    @{Row=35; Product='Bio Salmon'; Qty=1},
    @{Row=36; Product='Bio Chicken with Scallops'; Qty=1},
    @{Row=37; Product='whiskas pouch tuna-80 gr'; Qty=1},
    @{Row=38; Product='Meo pouch Otak² / Pempek'; Qty=1},
    @{Row=39; Product='Meo Pouch Tuna Kitten'; Qty=3},
    @{Row=40; Product='Pasir Top 5 L'; Qty=2},
    @{Row=41; Product='Furlove Tuna Kitten - 1kg'; Qty=1},
    @{Row=42; Product='Bolt Salmon - 800g'; Qty=2},
    @{Row=43; Product='Captain Cat Salmon Chicken'; Qty=1},
    @{Row=44; Product='Felibite Mom Kid Ikan'; Qty=1}
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Product
    $ws.Cells.Item($item.Row, 4).Value = $item.Qty
}

# Row 2 product cell gets a distinct font (Arial, black) per the new formatting
$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.Color = 0

# --- Clear old trailing rows 45-52 that no longer hold transaction data ---
$ws.Range("A45:F52").ClearContents()

# --- Extend the date-formatted column down through row 62 (new blank rows) ---
$ws.Range("A53:A62").NumberFormat = "yyyy\-mm\-dd;@"

# --- Update the visible selection to match the latest edit location ---
$ws.Range("F39:F44").Select()

